$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.288.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.01%  "

$ws.Range("D3").Value = "'3.307.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.09%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'191.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.17%  "

$ws.Range("D6").Value = "'559.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.47%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -2.37%  "

$ws.Range("D9").Value = "'3.298.85"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.13%  "

$ws.Range("E10").Value = "  -2.39%  "

$ws.Range("D11").Value = "'0.589"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.61%  "

$ws.Range("D12").Value = "'47.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.75%  "

$ws.Range("D13").Value = "'0.0000271"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.44%  "

$ws.Range("D14").Value = "'8.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.22%  "

$ws.Range("D15").Value = "'3.837.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.88%  "

$ws.Range("D16").Value = "'613.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.98%  "

$ws.Range("D17").Value = "'66.308.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").Value = "'18.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("E19").Value = "  -0.54%  "

$ws.Range("D20").Value = "'3.305.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.08%  "

$ws.Range("D21").Value = "'11.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.58%  "

$ws.Range("D22").Value = "'0.912"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.75%  "

$ws.Range("D23").Value = "'18.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.05%  "

$ws.Range("D24").Value = "'102.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.23%  "

$ws.Range("D25").Value = "'5.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.71%  "

$ws.Range("D26").Value = "'3.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.94%  "

$ws.Range("D27").Value = "'6.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.42%  "

$ws.Range("D28").Value = "'2.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.26%  "

$ws.Range("D29").Value = "'9.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.49%  "

$ws.Range("D30").Value = "'8.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.19%  "

$ws.Range("D31").Value = "'30.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.21%  "

$ws.Range("D32").Value = "'4.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.75%  "

$ws.Range("D33").Value = "'6.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.02%  "

$ws.Range("D34").Value = "'564.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.96%  "

$ws.Range("D35").Value = "'11.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.37%  "

$ws.Range("E36").Value = "  -0.81%  "

$ws.Range("D37").Value = "'3.763.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.74%  "

$ws.Range("E38").Value = "  -2.04%  "

$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("D40").Value = "'34.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.27%  "

$ws.Range("D41").Value = "'3.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.52%  "

$ws.Range("D42").Value = "'0.0₃0726"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.39%  "

$ws.Range("D43").Value = "'2.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("E44").Value = "  +0.51%  "

$ws.Range("D45").Value = "'0.341"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.29%  "

$ws.Range("D46").Value = "'3.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.24%  "

$ws.Range("E47").Value = "  +1.23%  "

$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("E49").Value = "  -1.35%  "

$ws.Range("E50").Value = "  -3.77%  "

$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.14%  "
